$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 600.2143
$ws.Range("I2").Value = 293.5
$ws.Range("J2").Value = 1009.1667
$ws.Range("K2").Value = 293.5
$ws.Range("L2").Value = 1009.1667
$ws.Range("M2").Value = -180.5
$ws.Range("N2").Value = -1235.1667
$ws.Range("H18").Value = 711.25
$ws.Range("I18").Value = 711.25
$ws.Range("K18").Value = 711.25
$ws.Range("M18").Value = -427.25
$ws.Range("H38").Value = 1602.1428
$ws.Range("J38").Value = 5933.3335
$ws.Range("L38").Value = 17800.0005
$ws.Range("N38").Value = -18544.0005
$ws.Range("H62").Value = 4900
$ws.Range("H65").Value = 4900
$ws.Range("H98").Value = 1594.8518
$ws.Range("I98").Value = 1436.409
$ws.Range("J98").Value = 2292
$ws.Range("K98").Value = 1436.409
$ws.Range("L98").Value = 2292
$ws.Range("M98").Value = 61.59099999999989
$ws.Range("N98").Value = -5288
$ws.Range("H112").Value = 911935.0600000001
$ws.Range("J112").Value = 1114176.4
$ws.Range("L112").Value = 3342529.2
$ws.Range("N112").Value = -3344745.2
$ws.Range("H122").Value = 1594.8518
$ws.Range("I122").Value = 1436.409
$ws.Range("J122").Value = 2292
$ws.Range("K122").Value = 4309.227000000001
$ws.Range("L122").Value = 6876
$ws.Range("M122").Value = -1859.227000000001
$ws.Range("N122").Value = -11776
$ws.Range("H137").Value = 8864
$ws.Range("I137").Value = 3548.3333
$ws.Range("K137").Value = 10644.9999
$ws.Range("M137").Value = -8094.999899999999
$ws.Range("H138").Value = 3561.13
$ws.Range("J138").Value = 3823.9102
$ws.Range("L138").Value = 11471.7306
$ws.Range("N138").Value = -21751.7306
$ws.Range("H141").Value = 5035.952
$ws.Range("I141").Value = 1710.3636
$ws.Range("K141").Value = 5131.0908
$ws.Range("M141").Value = 48.90920000000006

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1021.5263
$ws.Range("I2").Value = 1008.25714
$ws.Range("K2").Value = 1008.25714
$ws.Range("M2").Value = -895.25714
$ws.Range("H61").Value = 5199.2104
$ws.Range("I61").Value = 5159.9243
$ws.Range("J61").Value = 5719.75
$ws.Range("K61").Value = 5159.9243
$ws.Range("L61").Value = 5719.75
$ws.Range("M61").Value = -4947.9243
$ws.Range("N61").Value = -6143.75
$ws.Range("H63").Value = 2865
$ws.Range("I63").Value = 2865
$ws.Range("K63").Value = 2865
$ws.Range("M63").Value = -2179
$ws.Range("H66").Value = 2865
$ws.Range("I66").Value = 2865
$ws.Range("K66").Value = 14325
$ws.Range("M66").Value = -10893
$ws.Range("H74").Value = 3071.8235
$ws.Range("I74").Value = 1845.5217
$ws.Range("J74").Value = 5635.909
$ws.Range("K74").Value = 1845.5217
$ws.Range("L74").Value = 5635.909
$ws.Range("M74").Value = -971.5217
$ws.Range("N74").Value = -7383.909
$ws.Range("H77").Value = 3071.8235
$ws.Range("I77").Value = 1845.5217
$ws.Range("J77").Value = 5635.909
$ws.Range("K77").Value = 9227.6085
$ws.Range("L77").Value = 28179.545
$ws.Range("M77").Value = -4859.6085
$ws.Range("N77").Value = -36915.545
$ws.Range("H88").Value = 1942.3
$ws.Range("I88").Value = 2333
$ws.Range("J88").Value = 1774.8572
$ws.Range("K88").Value = 2333
$ws.Range("L88").Value = 1774.8572
$ws.Range("M88").Value = -1927
$ws.Range("N88").Value = -2586.8572
$ws.Range("H91").Value = 1942.3
$ws.Range("I91").Value = 2333
$ws.Range("J91").Value = 1774.8572
$ws.Range("K91").Value = 2333
$ws.Range("L91").Value = 1774.8572
$ws.Range("M91").Value = -929
$ws.Range("N91").Value = -4582.8572
$ws.Range("H116").Value = 1021.5263
$ws.Range("I116").Value = 1008.25714
$ws.Range("K116").Value = 1008.25714
$ws.Range("M116").Value = 1285.74286
$ws.Range("H125").Value = 110370
$ws.Range("J125").Value = 110370
$ws.Range("L125").Value = 110370
$ws.Range("N125").Value = -120210
$ws.Range("H132").Value = 2055.05
$ws.Range("I132").Value = 1801.5186
$ws.Range("K132").Value = 5404.5558
$ws.Range("M132").Value = -2874.5558
$ws.Range("H136").Value = 5199.2104
$ws.Range("I136").Value = 5159.9243
$ws.Range("J136").Value = 5719.75
$ws.Range("K136").Value = 15479.7729
$ws.Range("L136").Value = 17159.25
$ws.Range("M136").Value = -12929.7729
$ws.Range("N136").Value = -22259.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1021.5263
$ws.Range("I3").Value = 1008.25714
$ws.Range("K3").Value = 1008.25714
$ws.Range("M3").Value = -894.25714
$ws.Range("H86").Value = 3934057
$ws.Range("I86").Value = 6074115
$ws.Range("K86").Value = 6074115
$ws.Range("M86").Value = -6072992
$ws.Range("H89").Value = 3934057
$ws.Range("I89").Value = 6074115
$ws.Range("K89").Value = 30370575
$ws.Range("M89").Value = -30364959

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2309.15
$ws.Range("I16").Value = 1909.6154
$ws.Range("K16").Value = 1909.6154
$ws.Range("M16").Value = -1622.6154
$ws.Range("H31").Value = 3094.1724
$ws.Range("I31").Value = 2459.3333
$ws.Range("K31").Value = 2459.3333
$ws.Range("M31").Value = -2164.3333
$ws.Range("H34").Value = 3094.1724
$ws.Range("I34").Value = 2459.3333
$ws.Range("K34").Value = 2459.3333
$ws.Range("M34").Value = -2257.3333
$ws.Range("H58").Value = 1209.4546
$ws.Range("I58").Value = 1004.3461
$ws.Range("K58").Value = 1004.3461
$ws.Range("M58").Value = -801.3461
$ws.Range("H105").Value = 804.5714
$ws.Range("I105").Value = 852
$ws.Range("K105").Value = 852
$ws.Range("M105").Value = 895
$ws.Range("H113").Value = 2309.15
$ws.Range("I113").Value = 1909.6154
$ws.Range("K113").Value = 1909.6154
$ws.Range("M113").Value = 260.3846000000001
$ws.Range("H132").Value = 1334749.1
$ws.Range("I132").Value = 1334749.1
$ws.Range("K132").Value = 4004247.3
$ws.Range("M132").Value = -4001717.3
$ws.Range("H134").Value = 2264.4902
$ws.Range("I134").Value = 787.25
$ws.Range("K134").Value = 2361.75
$ws.Range("M134").Value = 173.25
$ws.Range("H136").Value = 1209.4546
$ws.Range("I136").Value = 1004.3461
$ws.Range("K136").Value = 3013.0383
$ws.Range("M136").Value = -463.0383000000002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 5972.75
$ws.Range("J113").Value = 5972.75
$ws.Range("L113").Value = 17918.25
$ws.Range("N113").Value = -22258.25
$ws.Range("H122").Value = 1382.2858
$ws.Range("J122").Value = 1409
$ws.Range("L122").Value = 12681
$ws.Range("N122").Value = -17581
$ws.Range("H138").Value = 9173000
$ws.Range("J138").Value = 10004500
$ws.Range("L138").Value = 30013500
$ws.Range("N138").Value = -30023780

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2876.9
$ws.Range("I126").Value = 2527.125
$ws.Range("K126").Value = 7581.375
$ws.Range("M126").Value = -5111.375
$ws.Range("H132").Value = 2728.7568
$ws.Range("I132").Value = 2728.7568
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8186.2704
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5656.2704
$ws.Range("N132").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2667.3809
$ws.Range("I22").Value = 1332.4286
$ws.Range("K22").Value = 1332.4286
$ws.Range("M22").Value = -1037.4286
$ws.Range("H27").Value = 2667.3809
$ws.Range("I27").Value = 1332.4286
$ws.Range("K27").Value = 1332.4286
$ws.Range("M27").Value = -1225.4286
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("I55").Value = 354.73334
$ws.Range("J55").Value = 332.7143
$ws.Range("K55").Value = 354.73334
$ws.Range("L55").Value = 332.7143
$ws.Range("M55").Value = -181.73334
$ws.Range("N55").Value = -678.7143
$ws.Range("H122").Value = 5153
$ws.Range("I122").Value = 4546.375
$ws.Range("K122").Value = 13639.125
$ws.Range("M122").Value = -11189.125
$ws.Range("H132").Value = 3435.0264
$ws.Range("I132").Value = 3449.6
$ws.Range("K132").Value = 10348.8
$ws.Range("M132").Value = -7818.799999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 62983.723
$ws.Range("J81").Value = 9576.200000000001
$ws.Range("L81").Value = 19152.4
$ws.Range("N81").Value = -21274.4
$ws.Range("H84").Value = 62983.723
$ws.Range("J84").Value = 9576.200000000001
$ws.Range("L84").Value = 95762
$ws.Range("N84").Value = -106370
$ws.Range("H113").Value = 3473068.2
$ws.Range("I113").Value = 7575890
$ws.Range("J113").Value = 1450.0769
$ws.Range("K113").Value = 22727670
$ws.Range("L113").Value = 4350.2307
$ws.Range("M113").Value = -22725500
$ws.Range("N113").Value = -8690.2307
$ws.Range("H132").Value = 1457.3462
$ws.Range("I132").Value = 1574.4375
$ws.Range("K132").Value = 4723.3125
$ws.Range("M132").Value = -2193.3125
$ws.Range("H136").Value = 31499.285
$ws.Range("I136").Value = 42099.6
$ws.Range("K136").Value = 126298.8
$ws.Range("M136").Value = -123748.8
